$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A18").Value = ""
$ws.Range("A18").Interior.Color = 65535
Write-Host "done"
